$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.722.60"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "1.630.49"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'214.39"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("E6").Value = "  -0.78%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  -1.30%  "
$ws.Range("D9").Value = "'0.0632"
$ws.Range("E9").Value = "  -1.41%  "
$ws.Range("D10").Value = "'19.48"
$ws.Range("E10").Value = "  -1.48%  "
$ws.Range("E11").Value = "  +1.48%  "
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").Value = "1.855.77"
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("D14").Value = "1.610.02"
$ws.Range("E14").Value = "  -3.55%  "
$ws.Range("D15").Value = "'0.555"
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").Value = "0.0₃0760"
$ws.Range("E16").Value = "  -2.15%  "
$ws.Range("D17").Value = "'62.84"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("D18").Value = "25.728.62"
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").Value = "'191.96"
$ws.Range("E21").Value = "  -1.14%  "
$ws.Range("D22").Value = "'9.92"
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("E23").Value = "  +1.44%  "
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("E25").Value = "  +3.16%  "
$ws.Range("D26").Value = "'142.27"
$ws.Range("E26").Value = "  +1.77%  "
$ws.Range("E27").Value = "  +1.79%  "
$ws.Range("D28").Value = "'6.85"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "'15.45"
$ws.Range("E29").Value = "  -1.13%  "
$ws.Range("D30").Value = "'1.24"
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("E33").Value = "  -1.20%  "
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("D36").Value = "'0.904"
$ws.Range("E36").Value = "  +0.74%  "
$ws.Range("D37").Value = "1.140.05"
$ws.Range("E37").Value = "  +2.95%  "
$ws.Range("E38").Value = "  -2.62%  "
$ws.Range("D39").Value = "'0.542"
$ws.Range("E39").Value = "  -2.07%  "
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("D41").Value = "'2.54"
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("D43").Value = "'5.54"
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("D44").Value = "'100.50"
$ws.Range("E44").Value = "  +1.29%  "
$ws.Range("D45").Value = "'0.804"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "1.765.77"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").Value = "0.0₆0108"
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("D48").Value = "'55.20"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").Value = "'0.0508"
$ws.Range("E49").Value = "  +0.84%  "
$ws.Range("D50").Value = "'0.419"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("D51").Value = "'1.43"
$ws.Range("E51").Value = "  +4.25%  "
